# Update odds values in Sheet1 to match the latest Betfair Back/Lay snapshot
# for 2026-01-07 (Jogos_do_Dia_Betfair_Back_Lay_2026-01-07.xlsx).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Napoli x Verona)
$ws.Range("I2").Value = 12.5
$ws.Range("J2").Value = 5.3
$ws.Range("K2").Value = 5.4
$ws.Range("Q2").Value = 2.1
$ws.Range("R2").Value = 1.32
$ws.Range("T2").Value = 2.6
$ws.Range("V2").Value = 1.08
$ws.Range("W2").Value = 3.6
$ws.Range("AB2").Value = 6.4
$ws.Range("AH2").Value = 40
$ws.Range("AI2").Value = 250
$ws.Range("AM2").Value = 330
$ws.Range("AN2").Value = 7.8

# Row 3 (Bologna x Atalanta)
$ws.Range("M3").Value = 1.07
$ws.Range("AI3").Value = 44

# Row 4 (Man City x Brighton)
$ws.Range("M4").Value = 1.02
$ws.Range("N4").Value = 6.8
$ws.Range("R4").Value = 1.81
$ws.Range("S4").Value = 2.2
$ws.Range("U4").Value = 2.4
$ws.Range("X4").Value = 36
$ws.Range("AB4").Value = 13

# Row 6 (Everton x Wolves)
$ws.Range("G6").Value = 1.83
$ws.Range("AF6").Value = 9.800000000000001
$ws.Range("AJ6").Value = 18.5
$ws.Range("AN6").Value = 14.5

# Row 7 (Crystal Palace x Aston Villa)
$ws.Range("I7").Value = 2.44
$ws.Range("N7").Value = 3.9
$ws.Range("Q7").Value = 1.98
$ws.Range("V7").Value = 1.69
$ws.Range("AG7").Value = 14

# Row 8 (Brentford x Sunderland)
$ws.Range("N8").Value = 3.6
$ws.Range("O8").Value = 1.36
$ws.Range("P8").Value = 1.87
$ws.Range("R8").Value = 1.33
$ws.Range("X8").Value = 13
$ws.Range("Y8").Value = 15
$ws.Range("AF8").Value = 10.5
$ws.Range("AM8").Value = 120

# Row 9 (Fulham x Chelsea)
$ws.Range("P9").Value = 2.14

# Row 10 (Lazio x Fiorentina)
$ws.Range("F10").Value = 2.54
$ws.Range("H10").Value = 3.35
$ws.Range("I10").Value = 3.4
$ws.Range("O10").Value = 1.46
$ws.Range("P10").Value = 1.68
$ws.Range("T10").Value = 2
$ws.Range("U10").Value = 1.94
$ws.Range("V10").Value = 1.41
$ws.Range("X10").Value = 9.6

# Row 11 (Torino x Udinese)
$ws.Range("L11").Value = 1.53
$ws.Range("Q11").Value = 2.44

# Row 12 (Parma x Inter)
$ws.Range("J12").Value = 5.4
$ws.Range("Q12").Value = 1.68
$ws.Range("T12").Value = 1.99
$ws.Range("U12").Value = 1.95
$ws.Range("AO12").Value = 5.6

# Row 13 (Newcastle x Leeds)
$ws.Range("K13").Value = 4.2

# Row 14 (Burnley x Man Utd)
$ws.Range("H14").Value = 1.72
$ws.Range("I14").Value = 1.73
$ws.Range("J14").Value = 4.2
$ws.Range("K14").Value = 4.3
$ws.Range("N14").Value = 4.5
$ws.Range("T14").Value = 1.8
$ws.Range("V14").Value = 2.36
$ws.Range("Y14").Value = 9.6
$ws.Range("AC14").Value = 9
